$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was re-sorted by "Casos totales" (col B) after new figures came
# in for several countries; this shuffles country rows (reordering sharedStrings)
# and updates the daily case/death counters for both the moved rows and some
# rows whose rank did not change. Write out the final state of every touched row.

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 4098478
$ws.Range("C4").Value = 69909
$ws.Range("D4").Value = 1938714
$ws.Range("E4").Value = 2013657
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1154
$ws.Range("H4").Value = 146107

# Row 5: Brasil -> Brasil
$ws.Range("A5").Value = 'Brasil'
$ws.Range("B5").Value = 2231871
$ws.Range("C5").Value = 65339
$ws.Range("D5").Value = 1532138
$ws.Range("E5").Value = 616843
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1293
$ws.Range("H5").Value = 82890

# Row 9: Peru -> Peru
$ws.Range("A9").Value = 'Peru'
$ws.Range("B9").Value = 366550
$ws.Range("C9").Value = 4463
$ws.Range("D9").Value = 252246
$ws.Range("E9").Value = 100725
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 13579

# Row 23: Argentina -> Argentina
$ws.Range("A23").Value = 'Argentina'
$ws.Range("B23").Value = 141900
$ws.Range("C23").Value = 5782
$ws.Range("D23").Value = 60531
$ws.Range("E23").Value = 78781
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 98
$ws.Range("H23").Value = 2588

# Row 59: Irlanda -> Japon
$ws.Range("A59").Value = 'Japon'
$ws.Range("B59").Value = 26303
$ws.Range("C59").Value = 567
$ws.Range("D59").Value = 20651
$ws.Range("E59").Value = 4663
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 989

# Row 60: Japon -> Irlanda
$ws.Range("A60").Value = 'Irlanda'
$ws.Range("B60").Value = 25819
$ws.Range("C60").Value = 17
$ws.Range("D60").Value = 23364
$ws.Range("E60").Value = 701
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 1754

# Row 70: Costa de Marfil -> Costa de Marfil
$ws.Range("A70").Value = 'Costa de Marfil'
$ws.Range("B70").Value = 14733
$ws.Range("C70").Value = 202
$ws.Range("D70").Value = 8995
$ws.Range("E70").Value = 5645
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 93

# Row 71: Chequia -> Chequia
$ws.Range("A71").Value = 'Chequia'
$ws.Range("B71").Value = 14570
$ws.Range("C71").Value = 246
$ws.Range("D71").Value = 9144
$ws.Range("E71").Value = 5062
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 364

# Row 79: Sudan -> Sudan
$ws.Range("A79").Value = 'Sudan'
$ws.Range("B79").Value = 11237
$ws.Range("C79").Value = 110
$ws.Range("D79").Value = 5835
$ws.Range("E79").Value = 4694
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 708

# Row 93: Guinea -> Guinea
$ws.Range("A93").Value = 'Guinea'
$ws.Range("B93").Value = 6747
$ws.Range("C93").Value = 95
$ws.Range("D93").Value = 5891
$ws.Range("E93").Value = 815
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 41

# Row 94: Gabon -> Gabon
$ws.Range("A94").Value = 'Gabon'
$ws.Range("B94").Value = 6588
$ws.Range("C94").Value = 155
$ws.Range("D94").Value = 4235
$ws.Range("E94").Value = 2306
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 47

# Row 95: Mauritania -> Mauritania
$ws.Range("A95").Value = 'Mauritania'
$ws.Range("B95").Value = 6027
$ws.Range("C95").Value = 42
$ws.Range("D95").Value = 3977
$ws.Range("E95").Value = 1895
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 155

# Row 98: Republica de Africa Central -> Republica de Africa Central
$ws.Range("A98").Value = 'Republica de Africa Central'
$ws.Range("B98").Value = 4574
$ws.Range("C98").Value = 13
$ws.Range("D98").Value = 1437
$ws.Range("E98").Value = 3080
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 57

# Row 103: Paraguay -> Paraguay
$ws.Range("A103").Value = 'Paraguay'
$ws.Range("B103").Value = 4000
$ws.Range("C103").Value = 183
$ws.Range("D103").Value = 2391
$ws.Range("E103").Value = 1573
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 36

# Row 107: Malaui -> Somalia
$ws.Range("A107").Value = 'Somalia'
$ws.Range("B107").Value = 3161
$ws.Range("C107").Value = 26
$ws.Range("D107").Value = 1495
$ws.Range("E107").Value = 1573
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 93

# Row 108: Somalia -> Malaui
$ws.Range("A108").Value = 'Malaui'
$ws.Range("B108").Value = 3149
$ws.Range("C108").Value = 104
$ws.Range("D108").Value = 1256
$ws.Range("E108").Value = 1822
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 7
$ws.Range("H108").Value = 71

# Row 131: Ruanda -> Ruanda
$ws.Range("A131").Value = 'Ruanda'
$ws.Range("B131").Value = 1689
$ws.Range("C131").Value = 34
$ws.Range("D131").Value = 867
$ws.Range("E131").Value = 817
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 5

# Row 139: Surinam -> Surinam
$ws.Range("A139").Value = 'Surinam'
$ws.Range("B139").Value = 1176
$ws.Range("C139").Value = 45
$ws.Range("D139").Value = 739
$ws.Range("E139").Value = 416
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 21

# Row 140: Jordania -> Niger
$ws.Range("A140").Value = 'Niger'
$ws.Range("B140").Value = 1122
$ws.Range("C140").Value = 9
$ws.Range("D140").Value = 1018
$ws.Range("E140").Value = 35
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 69

# Row 141: Liberia -> Jordania
$ws.Range("A141").Value = 'Jordania'
$ws.Range("B141").Value = 1120
$ws.Range("C141").Value = 7
$ws.Range("D141").Value = 1035
$ws.Range("E141").Value = 74
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 11

# Row 142: Niger -> Liberia
$ws.Range("A142").Value = 'Liberia'
$ws.Range("B142").Value = 1114
$ws.Range("C142").Value = 6
$ws.Range("D142").Value = 592
$ws.Range("E142").Value = 452
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 70

# Row 173: Islas Caimanes -> Bahamas
$ws.Range("A173").Value = 'Bahamas'
$ws.Range("B173").Value = 219
$ws.Range("C173").Value = 25
$ws.Range("D173").Value = 91
$ws.Range("E173").Value = 117
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 11

# Row 174: Camboya -> Islas Caimanes
$ws.Range("A174").Value = 'Islas Caimanes'
$ws.Range("B174").Value = 203
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 201
$ws.Range("E174").Value = 1
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 1

# Row 175: Guadalupe -> Camboya
$ws.Range("A175").Value = 'Camboya'
$ws.Range("B175").Value = 197
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 140
$ws.Range("E175").Value = 57
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

# Row 176: Bahamas -> Guadalupe
$ws.Range("A176").Value = 'Guadalupe'
$ws.Range("B176").Value = 195
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 172
$ws.Range("E176").Value = 9
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 14

# Row 181: Brunei -> Trinidad yTobago
$ws.Range("A181").Value = 'Trinidad yTobago'
$ws.Range("B181").Value = 141
$ws.Range("C181").Value = 2
$ws.Range("D181").Value = 127
$ws.Range("E181").Value = 6
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 8

# Row 182: Trinidad yTobago -> Brunei
$ws.Range("A182").Value = 'Brunei'
$ws.Range("B182").Value = 141
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 138
$ws.Range("E182").Value = 0
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 3

# Row 210: Islas Malvinas -> Groenlandia
$ws.Range("A210").Value = 'Groenlandia'
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Groenlandia -> Islas Malvinas
$ws.Range("A211").Value = 'Islas Malvinas'
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# --- Update the "last updated" timestamp footer (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 01:49"
